$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:K2").Value = -19.18438142348476
$ws.Range("C2").Value = 2.438596688683863

$ws.Range("B3:K3").Value = -19.18438142348476
$ws.Range("I3").Value = 2.439159412602895

$ws.Range("B4:K4").Value = -19.18438142348476
$ws.Range("C4").Value = 2.153082052256952
$ws.Range("D4").Value = 2.85470797420183
$ws.Range("F4").Value = 2.474149550625333
$ws.Range("H4").Value = 1.829577049245308
$ws.Range("J4").Value = 2.520694474768805

$ws.Range("B5:K5").Value = -19.18438142348476
$ws.Range("C5").Value = 0.9610893112447532
$ws.Range("G5").Value = 2.284434003888079

$ws.Range("B6:K6").Value = -19.18438142348476

$ws.Range("B7:K7").Value = -19.18438142348476
$ws.Range("B7").Value = 2.985986537256989

$ws.Range("B8:K8").Value = -19.18438142348476
$ws.Range("E8").Value = 2.916374902074046

$ws.Range("B9:K9").Value = -19.18438142348476
$ws.Range("B9").Value = 3.594226145620481

$ws.Range("B10:K10").Value = -19.18438142348476
$ws.Range("I10").Value = 1.61788012134753

$ws.Range("B11:K11").Value = -19.18438142348476
$ws.Range("E11").Value = 1.966130084749124
$ws.Range("G11").Value = 2.732862801945739

$ws.Range("B12:K12").Value = -19.18438142348476

$ws.Range("B13:K13").Value = -19.18438142348476
$ws.Range("E13").Value = 1.648543348860181
$ws.Range("J13").Value = 2.214623670629435

$ws.Range("B14:K14").Value = -19.18438142348476
$ws.Range("D14").Value = 1.650812403969608
$ws.Range("K14").Value = 4.321925794376789

$ws.Range("B15:K15").Value = -19.18438142348476
$ws.Range("D15").Value = -0.2800143451426065

$ws.Range("B16:K16").Value = -19.18438142348476
$ws.Range("J16").Value = 2.309339677435644

$ws.Range("B17:K17").Value = -19.18438142348476
$ws.Range("C17").Value = 0.6617532000434581
$ws.Range("D17").Value = -0.06003541502016801
$ws.Range("H17").Value = 0.5498535305201973
$ws.Range("I17").Value = 0.7238895093311771
$ws.Range("J17").Value = 1.166619247001618

$ws.Range("B18:K18").Value = -19.18438142348476
$ws.Range("H18").Value = 0.3722537261408301
$ws.Range("I18").Value = 0.8168082058472013
$ws.Range("J18").Value = 1.274038207080666

$ws.Range("B19:K19").Value = -19.18438142348476
$ws.Range("D19").Value = 1.718483109778024
$ws.Range("H19").Value = 2.027775911525903
$ws.Range("I19").Value = 2.093253897860568

$ws.Range("B20:K20").Value = -19.18438142348476
$ws.Range("C20").Value = 1.647770747322808
$ws.Range("D20").Value = 2.186683161071694
$ws.Range("F20").Value = 3.85235579943736
$ws.Range("H20").Value = 2.175833936196188
$ws.Range("I20").Value = 1.936310186769369

$ws.Range("B21:K21").Value = -19.18438142348476
$ws.Range("C21").Value = 1.794433057341916
$ws.Range("E21").Value = 2.435135132346153
$ws.Range("G21").Value = 3.084150582685556
$ws.Range("H21").Value = 2.347980685269695

